$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the 2021 column (R) of data, mirroring the structure of the existing
# 2020 column (Q) for every row of the table.
# ---------------------------------------------------------------------------

# Row 4 - year header (copy format from Q4, then set value)
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 2021

# Helper block data: for each "region" section we have 4 rows:
#   header (empty/percentage), share-of-damage row (formula), damage-volume
#   row (plain number) and GDP row (plain number, sometimes with heavier
#   border on the last one).

# ---- Кыргызская Республика (rows 5-8) -------------------------------------
$ws.Range("Q5").Copy($ws.Range("R5"))

$ws.Range("Q6").Copy($ws.Range("R6"))
$ws.Range("R6").Formula = "=R7/R8*100"

$ws.Range("Q7").Copy($ws.Range("R7"))
$ws.Range("R7").Value = 1931.83

$ws.Range("Q8").Copy($ws.Range("R8"))
$ws.Range("R8").Value = 739818.5

# ---- Баткенская область (rows 9-12) ----------------------------------------
$ws.Range("Q9").Copy($ws.Range("R9"))

$ws.Range("Q10").Copy($ws.Range("R10"))
$ws.Range("R10").Formula = "=R11/R12*100"

$ws.Range("Q11").Copy($ws.Range("R11"))
$ws.Range("R11").Value = 1552.9

$ws.Range("Q12").Copy($ws.Range("R12"))
$ws.Range("R12").Value = 25048.6

# ---- Джалал-Абадская область (rows 13-16) ----------------------------------
$ws.Range("Q13").Copy($ws.Range("R13"))

$ws.Range("Q14").Copy($ws.Range("R14"))
$ws.Range("R14").Formula = "=R15/R16*100"

$ws.Range("Q15").Copy($ws.Range("R15"))
$ws.Range("R15").Value = 125.7

$ws.Range("Q16").Copy($ws.Range("R16"))
$ws.Range("R16").Value = 82213.9
$ws.Range("R16").NumberFormat = "0.00"

# ---- Иссык-Кульская область (rows 17-20) -----------------------------------
$ws.Range("Q17").Copy($ws.Range("R17"))

$ws.Range("Q18").Copy($ws.Range("R18"))
$ws.Range("R18").Formula = "=R19/R20*100"

$ws.Range("Q19").Copy($ws.Range("R19"))
$ws.Range("R19").Value = 99.6

$ws.Range("Q20").Copy($ws.Range("R20"))
$ws.Range("R20").Value = 80059.6
$ws.Range("R20").NumberFormat = "0.00"

# ---- Нарынская область (rows 21-24) ----------------------------------------
$ws.Range("Q21").Copy($ws.Range("R21"))

$ws.Range("Q22").Copy($ws.Range("R22"))
$ws.Range("R22").Formula = "=R23/R24*100"

$ws.Range("Q23").Copy($ws.Range("R23"))
$ws.Range("R23").Value = 0.9

$ws.Range("Q24").Copy($ws.Range("R24"))
$ws.Range("R24").Value = 17172.7
$ws.Range("R24").NumberFormat = "0.00"

# ---- Ошская область (rows 25-28) -------------------------------------------
$ws.Range("Q25").Copy($ws.Range("R25"))

$ws.Range("Q26").Copy($ws.Range("R26"))
$ws.Range("R26").Formula = "=R27/R28*100"

$ws.Range("Q27").Copy($ws.Range("R27"))
$ws.Range("R27").Value = 15.9

$ws.Range("Q28").Copy($ws.Range("R28"))
$ws.Range("R28").Value = 56666.5
$ws.Range("R28").NumberFormat = "0.00"

# ---- Таласская область (rows 29-32) ----------------------------------------
$ws.Range("Q29").Copy($ws.Range("R29"))

$ws.Range("Q30").Copy($ws.Range("R30"))
$ws.Range("R30").Formula = "=R31/R32*100"

$ws.Range("Q31").Copy($ws.Range("R31"))
$ws.Range("R31").Value = 58.5

$ws.Range("Q32").Copy($ws.Range("R32"))
$ws.Range("R32").Value = 30765.1
$ws.Range("R32").NumberFormat = "0.00"

# ---- Чуйская область (rows 33-36) ------------------------------------------
$ws.Range("Q33").Copy($ws.Range("R33"))

$ws.Range("Q34").Copy($ws.Range("R34"))
$ws.Range("R34").Formula = "=R35/R36*100"

$ws.Range("Q35").Copy($ws.Range("R35"))
$ws.Range("R35").Value = 78.3

$ws.Range("Q36").Copy($ws.Range("R36"))
$ws.Range("R36").Value = 110267.1
$ws.Range("R36").NumberFormat = "0.00"

# ---- г. Бишкек (rows 37-40) -------------------------------------------------
$ws.Range("Q37").Copy($ws.Range("R37"))

$ws.Range("Q38").Copy($ws.Range("R38"))
$ws.Range("R38").Value = "-"

$ws.Range("Q39").Copy($ws.Range("R39"))
$ws.Range("R39").Value = "-"

$ws.Range("Q40").Copy($ws.Range("R40"))
$ws.Range("R40").Value = 297797.6
$ws.Range("R40").NumberFormat = "0.00"

# ---- г.Ош (rows 41-44) -------------------------------------------------------
$ws.Range("Q41").Copy($ws.Range("R41"))

$ws.Range("Q42").Copy($ws.Range("R42"))
$ws.Range("R42").Value = "-"

$ws.Range("Q43").Copy($ws.Range("R43"))
$ws.Range("R43").Value = "-"

$ws.Range("Q44").Copy($ws.Range("R44"))
$ws.Range("R44").Value = 39827.4
$ws.Range("R44").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# Match the saved selection state from the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("S14:T14").Select()
